$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (SYNTEX, 40x15)
$ws.Range("B2").Value = 5900
$ws.Range("C2").Value = 2950
$ws.Range("K2").Value = 2350

# Update row 3 (SYNTEX, 15x15)
$ws.Range("B3").Value = 5900
$ws.Range("C3").Value = 2950
$ws.Range("G3").Value = 15
$ws.Range("H3").Value = 15
$ws.Range("J3").Value = 225
$ws.Range("K3").Value = 2125
$ws.Range("L3").Value = 15
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 15

# Update row 4 (was INTRATEX 15x15, becomes SYNTEX 20x15 -- replaces old row4 content)
$ws.Range("A4").Value = "SYNTEX"
$ws.Range("B4").Value = 5900
$ws.Range("C4").Value = 2950
$ws.Range("D4").Value = "20x15"
$ws.Range("E4").Value = "(20, 25)"
$ws.Range("F4").Value = "arriba"
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 20
$ws.Range("I4").Value = 15
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 1825
$ws.Range("L4").Value = 20
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 20
$ws.Range("O4").Value = 15
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0

# Delete row 5 entirely (INTRATEX 20x15 row removed)
$ws.Rows(5).Delete()
